$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.468.36'
$ws.Range('E2').Value = '  -2.05%  '
$ws.Range('D3').Value = '3.475.33'
$ws.Range('E3').Value = '  -4.46%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '192.31'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.91%  '
$ws.Range('E7').Value = '  -2.75%  '
$ws.Range('D8').Value = '3.463.51'
$ws.Range('E8').Value = '  -4.51%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.204'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.617'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '51.43'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.00%  '
$ws.Range('E13').Value = '  -6.70%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.12'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.78%  '
$ws.Range('D15').Value = '4.032.94'
$ws.Range('E15').Value = '  -4.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '651.66'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.97%  '
$ws.Range('D17').Value = '69.349.05'
$ws.Range('E17').Value = '  -2.42%  '
$ws.Range('D18').Value = '3.483.58'
$ws.Range('E18').Value = '  -4.26%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.32'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.48%  '
$ws.Range('E20').Value = '  -1.84%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '18.18'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.85%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.945'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.02'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.41%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.27'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.87%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '99.06'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -6.85%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.28'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -7.65%  '
$ws.Range('E27').Value = '  -4.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.95'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.90%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.32'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.48'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.99%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.27'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -9.23%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.72'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -7.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.61'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.29%  '
$ws.Range('E34').Value = '  -5.76%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '60.86'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.22%  '
$ws.Range('D36').Value = '3.731.12'
$ws.Range('E36').Value = '  -6.31%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '528.54'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.85%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.998'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.03%  '
$ws.Range('D39').Value = '0.0₃0790'
$ws.Range('E39').Value = '  -9.66%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.91'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.50'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.27%  '
$ws.Range('E42').Value = '  -4.26%  '
$ws.Range('E43').Value = '  -1.90%  '
$ws.Range('B44').Value = 'CoreDAO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.51'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +72.64%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '34.30'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0441'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.34'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.82'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -9.48%  '
$ws.Range('E49').Value = '  -4.87%  '
$ws.Range('E50').Value = '  -0.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.16'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.31%  '
